$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timeline entry (row 37): Day 29, date 28/4/2024, 2 hours,
# "Made swagger with auth header, password reset with security code "
$newRow = $ws.Range("A37:D37")
$newRow.HorizontalAlignment = -4108   # xlCenter, matches the rest of the table
$newRow.VerticalAlignment = -4108    # xlCenter

$ws.Range("A37").Value = 29
$ws.Range("B37").Value = "28/4/2024"
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = "Made swagger with auth header, password reset with security code "

# Move the selection/view to the newly added row, like a user would after
# finishing data entry there.
$ws.Range("A37").Select() | Out-Null
